$d = $word.ActiveDocument

# Remove the stray "M2Doc version mismatch" warning block (its leading
# 4-space separator run and the orange "<---" marker run included),
# left over from a template/runtime version mismatch that no longer
# applies. The following "    " separator before "demonstration" is
# kept intact.
$r = $d.Content
$found = $r.Find.Execute(
    "    <---M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $r.Delete()
}
